# "Framework, Class, Basepage update"
#
# 1. The existing "Login" sheet gets a hyperlink on C2 (the text there,
#    "LmsHackathon@2024", looks like an email address and Excel's
#    autoformat turns it into a mailto: hyperlink + applies the built-in
#    "Hyperlink" cell style), and C2 becomes the active selection.
# 2. A new "Class" worksheet is added after "Login" and becomes the
#    active sheet, populated with a header row and a few sample rows.

$wb = $excel.ActiveWorkbook

# --- Login sheet: add the hyperlink on C2 ------------------------------
$loginWs = $wb.Worksheets.Item("Login")
$loginWs.Hyperlinks.Add($loginWs.Range("C2"), "mailto:LmsHackathon@2024") | Out-Null
$loginWs.Range("C2").Select() | Out-Null

# --- Add the new "Class" sheet right after "Login" ----------------------
$classWs = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Login"))
$classWs.Name = "Class"

# Header row
$headers = @("testcase", "Batchname", "ClassTopic", "ClassDescription", "SelectClassDates", "NoofClasses", "StaffName", "Status", "Comments", "Notes", "Recording")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $classWs.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Sample data rows (written in this order so the shared-string table ends
# up in the same order as the source workbook)
$classWs.Range("A2").Value = "validInputAll"
$classWs.Range("A3").Value = "validInputMandatory"
$classWs.Range("A4").Value = "Invalid"
$classWs.Range("C2").Value = "Oracle"

$classWs.Range("G9").Select() | Out-Null
